$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The sheet is protected; unprotect it so the cells below can be edited.
$ws.Unprotect()

# Update the confidentiality / "as of" date notice (A13): 2021-06-10 -> 2021-06-14
$ws.Range("A13").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-06-14 for illustrative purposes only and are subject to change."

# The multi-line text above can trigger an automatic row-height adjustment;
# restore row 13 to its normal (non-custom) auto height, matching the source file.
$ws.Rows(13).AutoFit()

# Update Weight (D) and Percent Change (E) values for rows 2-10
$ws.Range("D2").Value = 0.1337943701845034
$ws.Range("E2").Value = 0.007592315656274939

$ws.Range("D3").Value = 0.1084291316384225
$ws.Range("E3").Value = -0.004038866775930794

$ws.Range("D4").Value = 0.109761329878349
$ws.Range("E4").Value = 0.002068648482991042

$ws.Range("D5").Value = 0.1186805208577142
$ws.Range("E5").Value = -0.002326302729528607

$ws.Range("D6").Value = 0.1196803042200031
$ws.Range("E6").Value = 0.001638672675132957

$ws.Range("D7").Value = 0.148327643837175
$ws.Range("E7").Value = -0.007631456491391408

$ws.Range("D8").Value = 0.1303010785634669
$ws.Range("E8").Value = 0.003656821378340336

$ws.Range("D9").Value = 0.1310256208203658
$ws.Range("E9").Value = 0.001448754584395751

$ws.Range("D10").Value = 0.9999999999999999
$ws.Range("E10").Value = 0.0002593216854978042

# Re-apply sheet protection (original sheet was protected).
$ws.Protect()
